# Refactored Parser structure. Fixed problems with reading size of classes and methods.
#
# The parser that produces this workbook used to mis-measure the number of
# source lines for a class/method (it reported 0 for several members that
# clearly span more than zero lines). After the parser fix, the
# "Number of Lines" columns on the classNumberOfLines and
# methodNumberOfLines sheets need to be corrected to the real values.
#
# Values are written as text (leading "'" forces Excel to keep the
# numeric-looking values as text instead of silently re-typing the cell as
# a number), matching how the rest of that column is already stored.

$wb = $excel.ActiveWorkbook

function Set-LineCount {
    param($ws, $row, $col, $value)
    $ws.Cells.Item($row, $col).Value = "'" + $value
}

# --- classNumberOfLines: Number of Lines per class (column B) ---------
$classSheet = $wb.Worksheets.Item("classNumberOfLines")

Set-LineCount $classSheet 2 2 "66"   # domain.Order
Set-LineCount $classSheet 3 2 "7"    # domain.OrderStatus
Set-LineCount $classSheet 6 2 "3"    # domain.OrderSource

# --- methodNumberOfLines: Number of Lines per method (column C) -------
$methodSheet = $wb.Worksheets.Item("methodNumberOfLines")

# domain.Order plain getters/setters/toString/builder() -> 3 lines each
for ($r = 2; $r -le 17; $r++) {
    Set-LineCount $methodSheet $r 3 "3"
}

Set-LineCount $methodSheet 18 3 "2"   # Order()
Set-LineCount $methodSheet 19 3 "9"   # Order(Long, Long, Long, int, int, OrderStatus, OrderSource)

# domain.OrderStatus enum plumbing methods -> 1 line each
for ($r = 20; $r -le 23; $r++) {
    Set-LineCount $methodSheet $r 3 "1"
}

# domain.Order$OrderBuilder rows (24-33) are already correct - untouched

# domain.OrderSource enum plumbing methods -> 1 line each
for ($r = 34; $r -le 37; $r++) {
    Set-LineCount $methodSheet $r 3 "1"
}
